$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = "`n"

# G2: "Success" -> "SUCCESS"
$ws.Range("G2").Value = "SUCCESS"

# I3..I9: update day-of-week test step blocks (values + names changed, typo fixed: dayOfweek -> dayOfWeek)
function Build-Block([int]$num, [string]$name) {
    $lines = @(
        "Response_dayOfWeekNumber: $num",
        "DB_dayOfWeekNumber: $num",
        "Response_dayOfWeekFullName: $name",
        "DB_dayOfWeekFullName: $name",
        "Response_dayOfWeekShortName: $name",
        "DB_dayOfWeekShortName: $name",
        ""
    )
    return ($lines -join ($nl + $nl))
}

$ws.Range("I3").Value = Build-Block 7 "Test_G"
$ws.Range("I4").Value = Build-Block 5 "Test_E"
$ws.Range("I5").Value = Build-Block 2 "Test_B"
$ws.Range("I6").Value = Build-Block 4 "Test_D"
$ws.Range("I7").Value = Build-Block 1 "Test_A"
$ws.Range("I8").Value = Build-Block 3 "Test_C"
$ws.Range("I9").Value = Build-Block 6 "Test_F"

# Re-fit row heights so the COM engine doesn't leave stray explicit row
# heights behind after the multi-line text assignments above (the source
# workbook has no explicit row heights for these rows).
for ($r = 3; $r -le 9; $r++) {
    $ws.Rows.Item($r).AutoFit()
}

# G10: clear value (becomes an empty, styled cell), removing the reference to the "404" shared string
$ws.Range("G10").Value = ""

# I10: updated JSON error payload for the new swagger response format
$jsonLines = @(
    "",
    "{",
    "`t""meta"":",
    "`t{",
    "`t`t""version"":""1.0.0"",",
    "`t`t""transactionId"":""d273932d-b087-4d23-aa98-6b5294dbe073"",",
    "`t`t""timeStamp"":""2020-06-18T07:50:15.091+0000"",",
    "`t`t""statusCode"":""404"",",
    "`t`t""message"":",
    "`t`t{",
    "`t`t`t""status"":""ERROR"",",
    "`t`t`t""internalMessage"":""Exception"",",
    "`t`t`t""data"":",
    "`t`t`t{",
    "`t`t`t`t""errorMessage"":""Request Submitted With Error""",
    "`t`t`t}",
    "`t`t}",
    "`t},",
    "`t""errors"":",
    "`t[",
    "`t`t{",
    "`t`t`t""fieldName"":""Error"",",
    "`t`t`t""message"":""Could not find the GET method for URL /v2/daysOfWee""",
    "`t`t}",
    "`t]",
    "}"
)
$ws.Range("I10").Value = ($jsonLines -join $nl)
$ws.Rows.Item(10).AutoFit()
